# flash分配.xlsx -- "big improvement: add flash read/write"
#
# Adds 4 new "Light_threshold" (光敏门限) rows to the UPSSA0 block and a
# new "delay_time_num" (延时感应时间设置) row, replacing the old mostly
# empty separator rows. The existing resetbtcnt / bt_join_cnt rows shift
# down to make room.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header row (unchanged content, kept for completeness)
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "序号"
$ws.Range("B1").Value = "全局变量"
$ws.Range("C1").Value = "数据类型"
$ws.Range("D1").Value = "说明"
$ws.Range("E1").Value = "基地址"
$ws.Range("F1").Value = "偏移地址"
$ws.Range("G1").Value = "出厂值"

# ---------------------------------------------------------------------
# Rows 2-9: existing quick/slow detection fields (unchanged values)
# ---------------------------------------------------------------------
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "quick_time_times"
$ws.Range("C2").Value = "f32"
$ws.Range("D2").Value = "快检测时域乘法限"
$ws.Range("E2").Value = "UPSSA0"
$ws.Range("F2").Value = "0x0"
$ws.Range("G2").Value = 4

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "quick_time_add"
$ws.Range("C3").Value = "f32"
$ws.Range("D3").Value = "快检测时域加法限"
$ws.Range("E3").Value = "UPSSA0"
$ws.Range("F3").Value = "0x4"
$ws.Range("G3").Value = 32

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "quick_freq_times"
$ws.Range("C4").Value = "f32"
$ws.Range("D4").Value = "快检测频域乘法限"
$ws.Range("E4").Value = "UPSSA0"
$ws.Range("F4").Value = "0x8"
$ws.Range("G4").Value = 3

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "slow_time_times"
$ws.Range("C5").Value = "f32"
$ws.Range("D5").Value = "慢检测时域乘法限"
$ws.Range("E5").Value = "UPSSA0"
$ws.Range("F5").Value = "0xC"
$ws.Range("G5").Value = 4

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "slow_time_add"
$ws.Range("C6").Value = "f32"
$ws.Range("D6").Value = "慢检测时域加法限"
$ws.Range("E6").Value = "UPSSA0"
$ws.Range("F6").Value = "0x10"
$ws.Range("G6").Value = 40

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "slow_freq_times"
$ws.Range("C7").Value = "f32"
$ws.Range("D7").Value = "慢检测频域乘法限"
$ws.Range("E7").Value = "UPSSA0"
$ws.Range("F7").Value = "0x14"
$ws.Range("G7").Value = 6

$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "res_times"
$ws.Range("C8").Value = "f32"
$ws.Range("D8").Value = "慢检测呼吸限1"
$ws.Range("E8").Value = "UPSSA0"
$ws.Range("F8").Value = "0x18"
$ws.Range("G8").Value = 60

$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "offsetmin"
$ws.Range("C9").Value = "f32"
$ws.Range("D9").Value = "慢检测呼吸限2"
$ws.Range("E9").Value = "UPSSA0"
$ws.Range("F9").Value = "0x1C"
$ws.Range("G9").Value = 0.6

# ---------------------------------------------------------------------
# New rows 10-13: Light_threshold1..4 (光敏门限1..4) -- u32, UPSSA0 block
# Replaces the old near-empty row 10.
# ---------------------------------------------------------------------
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "Light_threshold1"
$ws.Range("C10").Value = "u32"
$ws.Range("D10").Value = "光敏门限1"
$ws.Range("E10").Value = "UPSSA0"
$ws.Range("F10").Value = "0x20"
$ws.Range("G10").Value = 0

$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "Light_threshold2"
$ws.Range("C11").Value = "u32"
$ws.Range("D11").Value = "光敏门限2"
$ws.Range("E11").Value = "UPSSA0"
$ws.Range("F11").Value = "0x24"
$ws.Range("G11").Value = 0

$ws.Range("A12").Value = 11
$ws.Range("B12").Value = "Light_threshold3"
$ws.Range("C12").Value = "u32"
$ws.Range("D12").Value = "光敏门限3"
$ws.Range("E12").Value = "UPSSA0"
$ws.Range("F12").Value = "0x28"
$ws.Range("G12").Value = 4000

$ws.Range("A13").Value = 12
$ws.Range("B13").Value = "Light_threshold4"
$ws.Range("C13").Value = "u32"
$ws.Range("D13").Value = "光敏门限4"
$ws.Range("E13").Value = "UPSSA0"
$ws.Range("F13").Value = "0x2C"
$ws.Range("G13").Value = 3800

# ---------------------------------------------------------------------
# Rows 14-15: resetbtcnt / bt_join_cnt, shifted down from old rows 11-12
# ---------------------------------------------------------------------
$ws.Range("A14").Value = 13
$ws.Range("B14").Value = "resetbtcnt"
$ws.Range("C14").Value = "u8"
$ws.Range("D14").Value = "蓝牙重新连接次数"
$ws.Range("E14").Value = "UPSSA1"
$ws.Range("F14").Value = "0x0"

$ws.Range("A15").Value = 14
$ws.Range("B15").Value = "bt_join_cnt"
$ws.Range("C15").Value = "u8"
$ws.Range("D15").Value = "蓝牙配网标志"
$ws.Range("E15").Value = "UPSSA1"
$ws.Range("F15").Value = "0x4"

# ---------------------------------------------------------------------
# New row 16: delay_time_num (延时感应时间设置) -- u32, replaces old
# MSA0 separator row 13.
# ---------------------------------------------------------------------
$ws.Range("A16").Value = 15
$ws.Range("B16").Value = "delay_time_num"
$ws.Range("C16").Value = "u32"
$ws.Range("D16").Value = "延时感应时间设置"
$ws.Range("E16").Value = "UPSSA0"
$ws.Range("F16").Value = "0x30"
$ws.Range("G16").Value = 0

# ---------------------------------------------------------------------
# Formatting: rows 10-16 need the same bordered look as the rest of the
# table. A2 already carries that exact "border + general number format"
# style, so copy it across the whole new block (this also switches the
# new G10:G16 cells away from the 0.000 numeric format used by G1:G9,
# onto a plain general one -- matching the new cellXfs entry added by
# the edit).
# ---------------------------------------------------------------------
$ws.Range("A2").Copy()
$ws.Range("A10:G16").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# Selection, as left by the editor after the edit
# ---------------------------------------------------------------------
$ws.Range("D21").Select()

Write-Host "flash-fenpei.xlsx updated: added Light_threshold1-4 and delay_time_num rows"
